$wb = $excel.ActiveWorkbook

# Add the new "short" worksheet after the existing "Blad1" sheet
$blad1 = $wb.Worksheets.Item("Blad1")
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $blad1)
$newSheet.Name = "short"

# Fill in header row
$newSheet.Range("A1").Value = "individual"
$newSheet.Range("B1").Value = "sex"
$newSheet.Range("C1").Value = "control"
$newSheet.Range("D1").Value = "treatment"

# Fill in data rows
$newSheet.Range("A2").Value = "A"
$newSheet.Range("B2").Value = "male"
$newSheet.Range("C2").Value = 23.6
$newSheet.Range("D2").Value = 21.3

$newSheet.Range("A3").Value = "B"
$newSheet.Range("B3").Value = "female"
$newSheet.Range("C3").Value = 19.3
$newSheet.Range("D3").Value = 17.9

$newSheet.Range("A4").Value = "C"
$newSheet.Range("B4").Value = "male"
$newSheet.Range("C4").Value = 25.8
$newSheet.Range("D4").Value = 24.1

# Select D1 on the new sheet and make it the active sheet/tab
$newSheet.Range("D1").Select()
$newSheet.Activate()
